# "retoques menu y slide" -- editorial.docx body copy edit:
#   - split the run " puedan salir adelante. ... los convoco a seguir trabajando "
#     into several runs, inserting "de " before "que el desarrollo...", turning
#     "del sector comercial y de las economias regionales" into "el sector
#     comercial, las economias regionales", and adding "la construccion, los
#     parques industriales, los jovenes empresarios y las mujeres que se dedican
#     a la actividad, " before "constituyen el programa...".
#   - Word's hidden _GoBack bookmark (marking the most recent edit point) moves
#     from the trailing empty paragraph to right after "regionales," as a side
#     effect of that edit, leaving the trailing paragraph bare.

$d = $word.ActiveDocument

# Locate the paragraph holding the sentence we need to edit rather than
# assuming a fixed paragraph index.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*puedan salir adelante*") {
        $targetPara = $cand
        break
    }
}
if ($targetPara -eq $null) {
    throw "could not locate the target paragraph"
}

# Rewrite the whole paragraph's WordOpenXML: every run keeps its original
# formatting (Verdana, sz 20), only the one big run gets split/edited, and the
# _GoBack bookmark is (re)inserted right after "regionales,".
$newParaXml = '<w:p w:rsidR="00C56CDE" w:rsidRDefault="00C56CDE" w:rsidP="00C56CDE"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Tenemos la certeza de que</w:t></w:r><w:r w:rsidR="00680536"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> la pandemia no acabará en 2021</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">, y </w:t></w:r><w:r w:rsidR="00680536"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">por ello </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>reivindicamos</w:t></w:r><w:r w:rsidR="00680536"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> el diálogo con todos los sectores y la</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> genuina demanda de asistencia para que nuestras pymes, que representa</w:t></w:r><w:r w:rsidR="00680536"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>n</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00680536"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>más del 40</w:t></w:r><w:r w:rsidR="00051649"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00680536"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>% del PBI nacional</w:t></w:r><w:r w:rsidR="00051649"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> y el 70 % del empleo registrado,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> puedan salir adelante. Porque apostamos por el país, como lo hemos hecho siempre, y lo seguiremos haciendo. Porque estamos convencidos </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">de </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">que el desarrollo de la pequeña y mediana industria, así como del </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>turismo, el sector comercial,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> las economías</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> regionales,</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> la construcción, los parques industriales, los jóvenes empresarios y las mujeres que se dedican a la actividad, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">constituyen el programa de crecimiento vertebral que la nación demanda en esta inédita etapa. Por todo ello, en honor a ese rol protagónico que nuestro sector representa, los convoco a seguir trabajando </w:t></w:r><w:r w:rsidR="004C013E"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">con valentía y convicción, </w:t></w:r><w:r w:rsidR="00680536"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">en unión, con </w:t></w:r><w:r w:rsidR="004C013E"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>entrega</w:t></w:r><w:r w:rsidR="00680536"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> y vocación, por un futuro mejor para todos.</w:t></w:r></w:p>'
$targetPara.Range.InsertXML($newParaXml) | Out-Null

# The trailing paragraph used to be the sole holder of the _GoBack bookmark;
# now that the bookmark lives earlier in the document, that paragraph goes
# back to being completely bare.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>") | Out-Null
